$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The MODS wrapper element changed from <update type="MODS"> to
# <datastream type="md_descriptive" operation="update">, so the shared
# strings that open/close that wrapper (stored in cells C2 and W2) need
# updating to match.
$ws.Range("C2").Value = '"><datastream type="md_descriptive" operation="update"><mods:mods xmlns:mods="http://www.loc.gov/mods/v3" xmlns:xlink="http://www.w3.org/1999/xlink" xmlns:xsi="http://www.w3.org/2001/XMLSchema-instance">'
$ws.Range("W2").Value = "</mods:mods></datastream></object>"

# Leave the selection on W2 to mirror the author's final cursor position.
$ws.Range("W2").Select()
